$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 5 ("Timeline & Milestones"): bold the "Phase 3" row in the table.
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$tbl5 = $s5.Shapes.Item(3).Table
for ($c = 1; $c -le $tbl5.Columns.Count; $c++) {
    $cell = $tbl5.Cell(4, $c)
    $cell.Shape.TextFrame.TextRange.Font.Bold = 1
}

# ---------------------------------------------------------------------------
# Slide 8 ("Investment Summary"): restructure the financial table.
# ---------------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$tbl8 = $s8.Shapes.Item(3).Table

# Insert two new columns after column 2 ("Year 1"):
#   -> AWS/Partner Credits (new col 3), Year 1 Net (new col 4)
$tbl8.Columns.Add(3)
$tbl8.Columns.Add(4)

# Append two new rows at the bottom for "Support & Maintenance" and
# "TOTAL INVESTMENT".
$tbl8.Rows.Add()
$tbl8.Rows.Add()

# Set the final column widths (EMU / 12700 = points).
$tbl8.Columns.Item(1).Width = 1567967 / 12700
$tbl8.Columns.Item(2).Width = 1132421 / 12700
$tbl8.Columns.Item(3).Width = 1742186 / 12700
$tbl8.Columns.Item(4).Width = 1132421 / 12700
$tbl8.Columns.Item(5).Width = 993046 / 12700
$tbl8.Columns.Item(6).Width = 993046 / 12700
$tbl8.Columns.Item(7).Width = 1158554 / 12700

function Set-HeaderCell($row, $col, $text) {
    $cell = $tbl8.Cell($row, $col)
    $cell.Shape.TextFrame.TextRange.Text = $text
    $cell.Shape.TextFrame.TextRange.Font.Size = 14
    $cell.Shape.TextFrame.TextRange.Font.Bold = 1
    $cell.Shape.TextFrame.TextRange.Font.Color.RGB = 0xFFFFFF
    $cell.Shape.Fill.Solid()
    $cell.Shape.Fill.ForeColor.RGB = 0x021CA0
}

function Set-DataCell($row, $col, $text, [bool]$bold) {
    $cell = $tbl8.Cell($row, $col)
    $cell.Shape.TextFrame.TextRange.Text = $text
    $cell.Shape.TextFrame.TextRange.Font.Size = 11
    if ($bold) {
        $cell.Shape.TextFrame.TextRange.Font.Bold = 1
    }
    $cell.Shape.Fill.Solid()
    $cell.Shape.Fill.ForeColor.RGB = 0xE6E6E7
}

# Row 1 (header)
Set-HeaderCell 1 1 'Cost Category'
Set-HeaderCell 1 2 'Year 1 List'
Set-HeaderCell 1 3 'AWS/Partner Credits'
Set-HeaderCell 1 4 'Year 1 Net'
Set-HeaderCell 1 5 'Year 2'
Set-HeaderCell 1 6 'Year 3'
Set-HeaderCell 1 7 '3-Year Total'

# Row 2: Professional Services
Set-DataCell 2 1 'Professional Services' $false
Set-DataCell 2 2 '$93,500' $false
Set-DataCell 2 3 '($8,000)' $false
Set-DataCell 2 4 '$85,500' $false
Set-DataCell 2 5 '$0' $false
Set-DataCell 2 6 '$0' $false
Set-DataCell 2 7 '$85,500' $false

# Row 3: Cloud Infrastructure
Set-DataCell 3 1 'Cloud Infrastructure' $false
Set-DataCell 3 2 '$8,644' $false
Set-DataCell 3 3 '($4,369)' $false
Set-DataCell 3 4 '$4,275' $false
Set-DataCell 3 5 '$8,644' $false
Set-DataCell 3 6 '$8,644' $false
Set-DataCell 3 7 '$21,563' $false

# Row 4: Software Licenses & Subscriptions
Set-DataCell 4 1 'Software Licenses & Subscriptions' $false
Set-DataCell 4 2 '$3,132' $false
Set-DataCell 4 3 '$0' $false
Set-DataCell 4 4 '$3,132' $false
Set-DataCell 4 5 '$3,132' $false
Set-DataCell 4 6 '$3,132' $false
Set-DataCell 4 7 '$9,396' $false

# Row 5: Support & Maintenance (new row)
Set-DataCell 5 1 'Support & Maintenance' $false
Set-DataCell 5 2 '$1,467' $false
Set-DataCell 5 3 '$0' $false
Set-DataCell 5 4 '$1,467' $false
Set-DataCell 5 5 '$1,467' $false
Set-DataCell 5 6 '$1,467' $false
Set-DataCell 5 7 '$4,401' $false

# Row 6: TOTAL INVESTMENT (new row, bold)
Set-DataCell 6 1 'TOTAL INVESTMENT' $true
Set-DataCell 6 2 '$106,743' $true
Set-DataCell 6 3 '($12,369)' $true
Set-DataCell 6 4 '$94,374' $true
Set-DataCell 6 5 '$13,243' $true
Set-DataCell 6 6 '$13,243' $true
Set-DataCell 6 7 '$120,860' $true
